{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list so it matches the\n// new, impact-focused accomplishment statements (see commit message /\n// diff). The section previously had 6 bullets; it now has 4. Several of\n// the old bullet strings (e.g. the trigonometric-algorithm line) also\n// appear verbatim elsewhere in the resume (the \"Partner - Siege Analytics\"\n// experience bullets), so we must scope our edits to only the paragraphs\n// that live inside the \"KEY ACHIEVEMENTS AND IMPACT\" section.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,styleBuiltIn\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the \"KEY ACHIEVEMENTS AND IMPACT\" Heading2 paragraph, and the next\n// Heading2 paragraph after it (the following section, \"TECHNICAL SKILLS\"),\n// to bound the section we are allowed to edit.\nconst sectionStart = items.findIndex(\n  (p) => p.styleBuiltIn === \"Heading2\" && p.text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\"\n);\nif (sectionStart === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\nlet sectionEnd = items.length;\nfor (let i = sectionStart + 1; i < items.length; i++) {\n  if (items[i].styleBuiltIn === \"Heading2\") {\n    sectionEnd = i;\n    break;\n  }\n}\nconst sectionItems = items.slice(sectionStart, sectionEnd);\n\n// Exact text we expect to find (old -> new), searched only within\n// sectionItems. The first three bullets are rewritten in place, and the\n// 4th rewritten bullet replaces what used to be the \"longitudinal data\n// analysis\" bullet (the last two old bullets, ETL pipelines / cloud data\n// warehouse, are removed outright).\nconst replacements = [\n  {\n    oldText: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    newText: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n  },\n  {\n    oldText: \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    newText: \"\u2022 $4.7M savings enabled nonprofit access\"\n  },\n  {\n    oldText: \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    newText: \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n  },\n  {\n    oldText: \"\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n    newText: \"\u2022 178% accuracy improvement in racial classification algorithms\"\n  }\n];\n\n// Paragraphs to remove outright (their text is no longer represented by any\n// bullet after the rewrite).\nconst toRemove = [\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n];\n\nconst usedIndices = new Set();\n\nfor (const { oldText, newText } of replacements) {\n  const idx = sectionItems.findIndex(\n    (p, i) => !usedIndices.has(i) && p.text.trim() === oldText\n  );\n  if (idx === -1) {\n    throw new Error(\"Could not find paragraph with text: \" + oldText);\n  }\n  usedIndices.add(idx);\n  sectionItems[idx].insertText(newText, \"Replace\");\n}\n\nconst removeIndices = [];\nfor (const target of toRemove) {\n  const idx = sectionItems.findIndex(\n    (p, i) => !usedIndices.has(i) && !removeIndices.includes(i) && p.text.trim() === target\n  );\n  if (idx === -1) {\n    throw new Error(\"Could not find paragraph to delete with text: \" + target);\n  }\n  removeIndices.push(idx);\n}\nfor (const idx of removeIndices) {\n  sectionItems[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list so it matches the\n# new, impact-focused accomplishment statements (see commit message /\n# diff). The section previously had 6 bullets; it now has 4. Several of\n# the old bullet strings (e.g. the trigonometric-algorithm line) also\n# appear verbatim elsewhere in the resume (the \"Partner - Siege Analytics\"\n# experience bullets), so we must scope our edits to only the paragraphs\n# that live inside the \"KEY ACHIEVEMENTS AND IMPACT\" section.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 paragraph.\n$sectionStartIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Heading 2\" -and $txt -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionStartIdx = $i\n        break\n    }\n}\nif ($sectionStartIdx -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Locate the next Heading 2 paragraph after it (bounds the section).\n$sectionEndIdx = $count + 1\nfor ($i = $sectionStartIdx + 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Heading 2\") {\n        $sectionEndIdx = $i\n        break\n    }\n}\n\n# Exact text we expect to find (old -> new), searched only within the\n# section's paragraph index range. The first three bullets are rewritten in\n# place, and the 4th rewritten bullet replaces what used to be the\n# \"longitudinal data analysis\" bullet (the last two old bullets, ETL\n# pipelines / cloud data warehouse, are removed outright).\n$oldText1 = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%'\n$newText1 = '\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%'\n\n$oldText2 = '\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations'\n$newText2 = '\u2022 $4.7M savings enabled nonprofit access'\n\n$oldText3 = '\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis'\n$newText3 = '\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions'\n\n$oldText4 = '\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality'\n$newText4 = '\u2022 178% accuracy improvement in racial classification algorithms'\n\n$oldText5 = '\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets'\n$oldText6 = '\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy'\n\n$oldTexts = @($oldText1, $oldText2, $oldText3, $oldText4)\n$newTexts = @($newText1, $newText2, $newText3, $newText4)\n$toRemove = @($oldText5, $oldText6)\n\n$matchedCount = 0\n$deleteIndices = @()\n\nfor ($i = $sectionStartIdx + 1; $i -lt $sectionEndIdx; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n\n    $matchIdx = -1\n    for ($k = 0; $k -lt $oldTexts.Count; $k++) {\n        if ($oldTexts[$k] -eq $txt) {\n            $matchIdx = $k\n            break\n        }\n    }\n\n    if ($matchIdx -ge 0) {\n        $p.Range.Text = $newTexts[$matchIdx]\n        $matchedCount = $matchedCount + 1\n    } elseif ($toRemove -contains $txt) {\n        $deleteIndices += $i\n    }\n}\n\nif ($matchedCount -ne 4) {\n    throw \"Expected to replace 4 paragraphs, replaced $matchedCount\"\n}\nif ($deleteIndices.Count -ne 2) {\n    throw \"Expected to find 2 paragraphs to delete, found $($deleteIndices.Count)\"\n}\n\n# Delete from highest index to lowest so earlier indices stay valid.\n$sortedDeleteIndices = $deleteIndices | Sort-Object -Descending\nforeach ($idx in $sortedDeleteIndices) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
